$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 498.91666
$ws.Range("I19").Value = 123.333336
$ws.Range("K19").Value = 123.333336
$ws.Range("M19").Value = 51.666664
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H103").Value = 999.75
$ws.Range("J103").Value = 1166.3334
$ws.Range("L103").Value = 3499.0002
$ws.Range("N103").Value = -4671.0002
$ws.Range("H138").Value = 3321.0715
$ws.Range("I138").Value = 623.75
$ws.Range("J138").Value = 4400
$ws.Range("K138").Value = 1871.25
$ws.Range("L138").Value = 13200
$ws.Range("M138").Value = 3268.75
$ws.Range("N138").Value = -23480

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 35043
$ws.Range("J53").Value = 35043
$ws.Range("L53").Value = 35043
$ws.Range("N53").Value = -36407
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H92").Value = 60000
$ws.Range("J92").Value = 60000
$ws.Range("L92").Value = 60000
$ws.Range("N92").Value = -64992
$ws.Range("H97").Value = 2925
$ws.Range("I97").Value = 2200
$ws.Range("K97").Value = 2200
$ws.Range("M97").Value = -1704
$ws.Range("H110").Value = 819
$ws.Range("I110").Value = 819
$ws.Range("K110").Value = 819
$ws.Range("M110").Value = 1226
$ws.Range("H132").Value = 4166.6665
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H36").Value = 7833
$ws.Range("I36").Value = 8499.5
$ws.Range("K36").Value = 8499.5
$ws.Range("M36").Value = -7965.5
$ws.Range("H86").Value = 2997.25
$ws.Range("I86").Value = 495
$ws.Range("K86").Value = 495
$ws.Range("M86").Value = 628
$ws.Range("H89").Value = 2997.25
$ws.Range("I89").Value = 495
$ws.Range("K89").Value = 2475
$ws.Range("M89").Value = 3141
$ws.Range("H105").Value = 1086.6666
$ws.Range("I105").Value = 880
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 880
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 867
$ws.Range("N105").Value = -4994
$ws.Range("H125").Value = 69333.336
$ws.Range("J125").Value = 69333.336
$ws.Range("L125").Value = 69333.336
$ws.Range("N125").Value = -79173.336
$ws.Range("H134").Value = 2365.5
$ws.Range("I134").Value = 2138.6
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 6415.799999999999
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -3880.799999999999
$ws.Range("N134").Value = -15570

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14337.125
$ws.Range("I31").Value = 5844.4287
$ws.Range("J31").Value = 20942.555
$ws.Range("K31").Value = 5844.4287
$ws.Range("L31").Value = 20942.555
$ws.Range("M31").Value = -5549.4287
$ws.Range("N31").Value = -21532.555
$ws.Range("H33").Value = 6423.75
$ws.Range("J33").Value = 17233.5
$ws.Range("L33").Value = 17233.5
$ws.Range("N33").Value = -17991.5
$ws.Range("H34").Value = 14337.125
$ws.Range("I34").Value = 5844.4287
$ws.Range("J34").Value = 20942.555
$ws.Range("K34").Value = 5844.4287
$ws.Range("L34").Value = 20942.555
$ws.Range("M34").Value = -5642.4287
$ws.Range("N34").Value = -21346.555
$ws.Range("H42").Value = 45000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H44").Value = 29999.5
$ws.Range("I44").Value = 29999
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 29999
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -29557
$ws.Range("N44").Value = -30884
$ws.Range("H62").Value = 3500
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876
$ws.Range("H65").Value = 3500
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380
$ws.Range("H99").Value = 718591.9
$ws.Range("I99").Value = 1714.3334
$ws.Range("J99").Value = 1256250
$ws.Range("K99").Value = 1714.3334
$ws.Range("L99").Value = 1256250
$ws.Range("M99").Value = -216.3334
$ws.Range("N99").Value = -1259246
$ws.Range("H105").Value = 419.5
$ws.Range("I105").Value = 419.5
$ws.Range("K105").Value = 419.5
$ws.Range("M105").Value = 1327.5
$ws.Range("H122").Value = 1002.2727
$ws.Range("J122").Value = 845.4
$ws.Range("L122").Value = 2536.2
$ws.Range("N122").Value = -7436.2
$ws.Range("H126").Value = 718591.9
$ws.Range("I126").Value = 1714.3334
$ws.Range("J126").Value = 1256250
$ws.Range("K126").Value = 5143.0002
$ws.Range("L126").Value = 3768750
$ws.Range("M126").Value = -2673.0002
$ws.Range("N126").Value = -3773690
$ws.Range("H132").Value = 1373
$ws.Range("I132").Value = 1268.5
$ws.Range("K132").Value = 3805.5
$ws.Range("M132").Value = -1275.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5687.857
$ws.Range("J3").Value = 4969.5
$ws.Range("L3").Value = 14908.5
$ws.Range("N3").Value = -15132.5
$ws.Range("H140").Value = 20000
$ws.Range("I140").Value = 20000
$ws.Range("K140").Value = 60000
$ws.Range("M140").Value = -54820
$ws.Range("H141").Value = 3030
$ws.Range("I141").Value = 3030
$ws.Range("K141").Value = 9090
$ws.Range("M141").Value = -3910

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 18000
$ws.Range("J58").Value = 18000
$ws.Range("L58").Value = 18000
$ws.Range("N58").Value = -18554
$ws.Range("H97").Value = 3157.8333
$ws.Range("I97").Value = 3157.8333
$ws.Range("K97").Value = 3157.8333
$ws.Range("M97").Value = -2661.8333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2916.3333
$ws.Range("J22").Value = 3324.75
$ws.Range("L22").Value = 3324.75
$ws.Range("N22").Value = -3914.75
$ws.Range("H26").Value = 3531.125
$ws.Range("I26").Value = 1249.5
$ws.Range("J26").Value = 4291.6665
$ws.Range("K26").Value = 1249.5
$ws.Range("L26").Value = 4291.6665
$ws.Range("M26").Value = -954.5
$ws.Range("N26").Value = -4881.6665
$ws.Range("H27").Value = 2916.3333
$ws.Range("J27").Value = 3324.75
$ws.Range("L27").Value = 3324.75
$ws.Range("N27").Value = -3538.75
$ws.Range("H39").Value = 600
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 419.42856
$ws.Range("I55").Value = 548.5
$ws.Range("J55").Value = 397.91666
$ws.Range("K55").Value = 548.5
$ws.Range("L55").Value = 397.91666
$ws.Range("M55").Value = -375.5
$ws.Range("N55").Value = -743.91666
$ws.Range("H56").Value = 45874.75
$ws.Range("J56").Value = 45833.332
$ws.Range("L56").Value = 45833.332
$ws.Range("N56").Value = -47215.332
$ws.Range("H58").Value = 46000
$ws.Range("J58").Value = 46000
$ws.Range("L58").Value = 46000
$ws.Range("N58").Value = -46520
$ws.Range("H122").Value = 996.5
$ws.Range("I122").Value = 996.5
$ws.Range("K122").Value = 2989.5
$ws.Range("M122").Value = -539.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 46333.332
$ws.Range("J52").Value = 39500
$ws.Range("L52").Value = 39500
$ws.Range("N52").Value = -39952
$ws.Range("H122").Value = 935.5714
$ws.Range("J122").Value = 1062.25
$ws.Range("L122").Value = 3186.75
$ws.Range("N122").Value = -8086.75
